$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190560698509216
$ws.Range("B1").Value = 2.206927537918091
$ws.Range("C1").Value = 10.55808639526367
$ws.Range("D1").Value = 2.5734543800354
$ws.Range("E1").Value = 1.231905817985535
